# Applies the task7 "http json mobile task" edit:
# - Updates the four text answers in column B (rows 2-5)
# - Updates the expense amounts in column B (rows 6, 10-15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text fields (become new shared strings: вчапсп / ваыа / аывпв / ыва - 43цкуы)
$ws.Range("B2").Value = "вчапсп"
$ws.Range("B3").Value = "ваыа"
$ws.Range("B4").Value = "аывпв"
$ws.Range("B5").Value = "ыва - 43цкуы"

# Numeric amounts
$ws.Range("B6").Value = 454
$ws.Range("B10").Value = 5454
$ws.Range("B11").Value = 534534
$ws.Range("B12").Value = 543543
$ws.Range("B13").Value = 1083531
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 1083077
